# Add data for 2022-01-02
# - Rename sheet (title reflects the "through" date moving from 12-24 to 12-25)
# - Update the December row label accordingly
# - Update December (row 13) and Total (row 14) figures with the new day's data

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename worksheet tab to reflect new "through" date
$ws.Name = "Through 2021-12-25"

# Update the December month label in column A (row 13)
$ws.Range("A13").Value = "December (through 12-25)"

# Update December row (row 13) values
$ws.Range("B13").Value = 36
$ws.Range("C13").Value = 82
$ws.Range("D13").Value = 97
$ws.Range("E13").Value = 57
$ws.Range("F13").Value = 52
$ws.Range("G13").Value = 118
$ws.Range("H13").Value = 158

# Update Total row (row 14) values
$ws.Range("B14").Value = 327
$ws.Range("C14").Value = 645
$ws.Range("D14").Value = 918
$ws.Range("E14").Value = 739
$ws.Range("F14").Value = 586
$ws.Range("G14").Value = 1382
$ws.Range("H14").Value = 1801
